$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.284.66"
$ws.Range("E2").Value = "  -3.57%  "
$ws.Range("D3").Value = "2.375.45"
$ws.Range("E3").Value = "  -3.66%  "
$ws.Range("E4").Value = "  -0.05%  "
$s = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "499.07"
$ws.Range("D5").Style = $s
$ws.Range("E5").Value = "  -5.68%  "
$s = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.06"
$ws.Range("D6").Style = $s
$ws.Range("E6").Value = "  -2.07%  "
$s = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = $s
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("D9").Value = "2.402.52"
$ws.Range("E9").Value = "  -2.55%  "
$s = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0956"
$ws.Range("D10").Style = $s
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("E11").Value = "  -0.93%  "
$s = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.63"
$ws.Range("D12").Style = $s
$ws.Range("E12").Value = "  -7.03%  "
$ws.Range("E13").Value = "  -2.09%  "
$ws.Range("D14").Value = "2.803.68"
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").Value = "56.168.43"
$ws.Range("E15").Value = "  -3.62%  "
$s = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.49"
$ws.Range("D16").Style = $s
$ws.Range("E16").Value = "  -2.13%  "
$ws.Range("E17").Value = "  -1.69%  "
$ws.Range("D18").Value = "2.390.46"
$ws.Range("E18").Value = "  -3.74%  "
$ws.Range("E19").Value = "  -3.46%  "
$s = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "309.33"
$ws.Range("D20").Style = $s
$ws.Range("E20").Value = "  -2.77%  "
$ws.Range("E21").Value = "  -3.91%  "
$s = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.25"
$ws.Range("D22").Style = $s
$ws.Range("E22").Value = "  +2.04%  "
$s = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = $s
$ws.Range("E23").Value = "  -0.18%  "
$s = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.58"
$ws.Range("D24").Style = $s
$ws.Range("E24").Value = "  -4.61%  "
$s = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "64.78"
$ws.Range("D25").Style = $s
$ws.Range("E25").Value = "  -0.65%  "
$s = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("D26").Style = $s
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("D27").Value = "2.488.19"
$ws.Range("E27").Value = "  -4.48%  "
$s = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.375"
$ws.Range("D28").Style = $s
$ws.Range("E28").Value = "  -6.56%  "
$s = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.149"
$ws.Range("D29").Style = $s
$ws.Range("E29").Value = "  -5.61%  "
$ws.Range("E30").Value = "  -0.44%  "
$s = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "172.54"
$ws.Range("D31").Style = $s
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").Value = "0.0₃0712"
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("E33").Value = "  -3.12%  "
$s = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.11"
$ws.Range("D34").Style = $s
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("E35").Value = "  -0.09%  "
$s = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.09"
$ws.Range("D36").Style = $s
$ws.Range("E36").Value = "  -6.45%  "
$ws.Range("E37").Value = "  -0.36%  "
$s = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.81"
$ws.Range("D38").Style = $s
$ws.Range("E38").Value = "  -0.30%  "
$s = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.21"
$ws.Range("D39").Style = $s
$ws.Range("E39").Value = "  +1.27%  "
$ws.Range("E40").Value = "  -1.50%  "
$s = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.80"
$ws.Range("D41").Style = $s
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$s = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.42"
$ws.Range("D42").Style = $s
$ws.Range("E42").Value = "  -3.80%  "
$ws.Range("B43").Value = "SuiNetwork"
$ws.Range("C43").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$s = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.785"
$ws.Range("D43").Style = $s
$ws.Range("E43").Value = "  -2.47%  "
$s = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "128.98"
$ws.Range("D44").Style = $s
$ws.Range("E44").Value = "  +0.47%  "
$s = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.34"
$ws.Range("D45").Style = $s
$ws.Range("E45").Value = "  -2.60%  "
$s = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.74"
$ws.Range("D46").Style = $s
$ws.Range("E46").Value = "  -3.23%  "
$s = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "251.72"
$ws.Range("D47").Style = $s
$ws.Range("E47").Value = "  -6.46%  "
$s = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.562"
$ws.Range("D48").Style = $s
$ws.Range("E48").Value = "  -3.71%  "
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("E50").Value = "  -2.91%  "
$s = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.82"
$ws.Range("D51").Style = $s
$ws.Range("E51").Value = "  -1.79%  "
